# Actualización desde MV -datos-
# Update the last existing data row (139) with revised figures and append
# a new monthly row (140) for 01-07-2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing row 139 values ---
$ws.Range("B139").Value = 125.5
$ws.Range("C139").Value = 97.1
$ws.Range("D139").Value = 125.9

# --- Append new row 140 ---
# Column A holds the period label as text (e.g. "01-06-2021"). Excel would
# normally auto-convert a dd-mm-yyyy-looking string into a date serial, so
# build it as a text formula first and then convert the formula result to a
# static value via copy / paste-special (values only). This keeps the cell
# stored as shared-string text (matching the existing column) instead of a
# date number or a live formula, and avoids introducing any new cell styles.
$ws.Range("A140").Formula = '="01-07-2021"'
$ws.Range("A140").Copy() | Out-Null
$ws.Range("A140").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B140").Value = 127.1
$ws.Range("C140").Value = 97.6
$ws.Range("D140").Value = 127.7
